$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.269.02'
$ws.Range("E2").Value = '  +1.44%  '

$ws.Range("D3").Value = '3.897.83'
$ws.Range("E3").Value = '  +0.17%  '

$ws.Range("E4").Value = '  +0.14%  '

$ws.Range("D5").Value = '''526.99'
$ws.Range("E5").Value = '  +8.97%  '

$ws.Range("D6").Value = '''144.75'
$ws.Range("E6").Value = '  -0.16%  '

$ws.Range("D7").Value = '''0.613'
$ws.Range("E7").Value = '  -1.57%  '

$ws.Range("E8").Value = '  +0.06%  '

$ws.Range("D9").Value = '''0.717'
$ws.Range("E9").Value = '  -3.28%  '

$ws.Range("E10").Value = '  -4.83%  '

$ws.Range("D11").Value = '''0.0000337'
$ws.Range("E11").Value = '  -4.43%  '

$ws.Range("D13").Value = '4.528.16'
$ws.Range("E13").Value = '  +0.42%  '

$ws.Range("D14").Value = '''10.26'
$ws.Range("E14").Value = '  -2.47%  '

$ws.Range("D15").Value = '3.899.49'
$ws.Range("E15").Value = '  -0.23%  '

$ws.Range("B16").Value = 'Uniswap'
$ws.Range("C16").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D16").Value = '''14.00'
$ws.Range("E16").Value = '  -1.29%  '

$ws.Range("B17").Value = 'Polygon'
$ws.Range("C17").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D17").Value = '''1.23'
$ws.Range("E17").Value = '  +8.49%  '

$ws.Range("E18").Value = '  -0.77%  '

$ws.Range("D19").Value = '''19.75'
$ws.Range("E19").Value = '  -0.90%  '

$ws.Range("D20").Value = '69.214.68'
$ws.Range("E20").Value = '  +1.40%  '

$ws.Range("D21").Value = '''426.18'
$ws.Range("E21").Value = '  -0.86%  '

$ws.Range("E22").Value = '  -6.62%  '

$ws.Range("B23").Value = 'InternetComputer(DFINITY)'
$ws.Range("C23").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D23").Value = '''14.14'
$ws.Range("E23").Value = '  -4.69%  '

$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = '''88.26'
$ws.Range("E24").Value = '  -0.98%  '

$ws.Range("D25").Value = '''4.03'
$ws.Range("E25").Value = '  +10.02%  '

$ws.Range("D26").Value = '''11.41'
$ws.Range("E26").Value = '  -8.60%  '

$ws.Range("D27").Value = '''10.59'
$ws.Range("E27").Value = '  -3.96%  '

$ws.Range("D28").Value = '''36.42'
$ws.Range("E28").Value = '  -2.38%  '

$ws.Range("D29").Value = '''678.89'
$ws.Range("E29").Value = '  -4.79%  '

$ws.Range("D30").Value = '''13.12'
$ws.Range("E30").Value = '  -2.77%  '

$ws.Range("E31").Value = '  -3.21%  '

$ws.Range("E32").Value = '  -3.09%  '

$ws.Range("D33").Value = '''68.62'
$ws.Range("E33").Value = '  +11.26%  '

$ws.Range("D34").Value = '0.0₃0878'
$ws.Range("E34").Value = '  +0.06%  '

$ws.Range("D35").Value = '''0.434'
$ws.Range("E35").Value = '  +8.64%  '

$ws.Range("E36").Value = '  -1.14%  '

$ws.Range("D37").Value = '''39.98'
$ws.Range("E37").Value = '  -2.50%  '

$ws.Range("D38").Value = '''0.148'
$ws.Range("E38").Value = '  +1.32%  '

$ws.Range("E39").Value = '  +0.02%  '

$ws.Range("E40").Value = '  -0.01%  '

$ws.Range("E41").Value = '  +5.85%  '

$ws.Range("D42").Value = '''0.0481'
$ws.Range("E42").Value = '  -3.50%  '

$ws.Range("D43").Value = '''3.18'
$ws.Range("E43").Value = '  +7.74%  '

$ws.Range("E44").Value = '  -8.36%  '

$ws.Range("E45").Value = '  -0.88%  '

$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").Value = '''0.140'
$ws.Range("E46").Value = '  -1.80%  '

$ws.Range("B47").Value = 'FLOKI'
$ws.Range("C47").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D47").Value = '''0.000284'
$ws.Range("E47").Value = '  +17.57%  '

$ws.Range("E48").Value = '  +6.77%  '

$ws.Range("D49").Value = '0.0₆0350'
$ws.Range("E49").Value = '  -0.04%  '

$ws.Range("D50").Value = '2.745.52'
$ws.Range("E50").Value = '  +14.15%  '

$ws.Range("D51").Value = '''145.30'
$ws.Range("E51").Value = '  +0.48%  '
